# Atualização automática de preços de eletricidade
# Updates row 2 of the Spot_PT sheet with the new day's hourly spot prices.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (date serial number, keeps existing date formatting)
$ws.Range("A2").Value = 45994

# Hourly prices 0h-1h .. 23h-24h (columns B..Y)
$ws.Range("B2").Value = 85.05
$ws.Range("C2").Value = 79.36
$ws.Range("D2").Value = 74.52
$ws.Range("E2").Value = 68.15000000000001
$ws.Range("F2").Value = 67.14
$ws.Range("G2").Value = 72.23
$ws.Range("H2").Value = 77.90000000000001
$ws.Range("I2").Value = 93.06999999999999
$ws.Range("J2").Value = 99.84
$ws.Range("K2").Value = 83.88
$ws.Range("L2").Value = 67.38
$ws.Range("M2").Value = 38.14
$ws.Range("N2").Value = 20.63
$ws.Range("O2").Value = 15.44
$ws.Range("P2").Value = 14.59
$ws.Range("Q2").Value = 31.11
$ws.Range("R2").Value = 64.47
$ws.Range("S2").Value = 88.29000000000001
$ws.Range("T2").Value = 102.44
$ws.Range("U2").Value = 99.22
$ws.Range("V2").Value = 96.88
$ws.Range("W2").Value = 88.76000000000001
$ws.Range("X2").Value = 85.48
$ws.Range("Y2").Value = 76.05

# Price_Daily_Avg
$ws.Range("Z2").Value = 70.42

# Slot_4h_max (unchanged) / Slot_4h_price
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 88.59999999999999

# Slot_2h_frist (unchanged) / Slot_2h_frist_price
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 100.83

# Slot_2h_second (unchanged) / Slot_2h_second_price
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 92.81999999999999

# Slot_min_price
$ws.Range("AG2").Value = "3h-16h"
